$wb = $excel.ActiveWorkbook

$wsExternal = $wb.Worksheets.Item("ExternalContact")
$wsEng = $wb.Worksheets.Item("EngContact")

# --- Update ExternalContact sheet content ---
# A1 previously held the header "External Contact Name" (bold). It now holds
# the value that used to be in A2 ("Aaron M. Rosen") with normal (non-bold) style.
$wsExternal.Range("A1").Value = "Aaron M. Rosen"
$wsExternal.Range("A1").Font.Bold = $false

# A2 now holds the new name "Ramana Sail"
$wsExternal.Range("A2").Value = "Ramana Sail"

# --- EngContact sheet content stays the same text-wise ---
$wsEng.Range("A1").Value = "Engagement Contact Name"
$wsEng.Range("A2").Value = "Aaron Rosen"

# --- Selections on each sheet ---
$wsEng.Range("C6").Select()
$wsExternal.Range("A2").Select()

# --- Activate ExternalContact as the active sheet/tab ---
$wsExternal.Activate()
